# Generate Report for Handback
# Updates the localization-status report to reflect a failed handback
# transform for the ac9ebc5a-bd66-427a-b324-9455da75c320 document, and
# records the error detail for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 3 corresponds to ac9ebc5a-bd66-427a-b324-9455da75c320.md
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = "Handback file name: bktzdym2.neq is different with handoff file name: ac9ebc5a-bd66-427a-b324-9455da75c320.55da7ee1e0bafe7c412adcc8a89d5ad91796aca7.zh-cn."
# ColumnWidth is stored internally with a small fixed padding offset relative
# to the saved OOXML "width" attribute, so back out that offset to land on 40.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet: row 3 corresponds to ac9ebc5a-bd66-427a-b324-9455da75c320.md
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = "Handback file name: bktzdym2.neq is different with handoff file name: ac9ebc5a-bd66-427a-b324-9455da75c320.55da7ee1e0bafe7c412adcc8a89d5ad91796aca7.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

# Overview sheet: the per-locale status columns mirror each sheet's Status
# value for this row, so they pick up the same new status text.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
